# Deploy the implementation guide.
#
# 1. Update the "Date" metadata value on the Metadata sheet.
# 2. Append a new concept row ("POLYM" / "Polymalformation") to the
#    Concepts sheet, reusing the existing row-10 formatting so no new
#    cell styles are introduced.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date value ---------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2023-02-14T19:15:12+00:00"

# --- Concepts sheet: add the new Polymalformation concept row ----------
$wsConcepts = $wb.Worksheets.Item("Concepts")

# Clone row 10 (values + style) into the new row 11 so the new row keeps
# the same cell style (s="2") used by every other data row.
$wsConcepts.Range("A10:D10").Copy($wsConcepts.Range("A11:D11"))

# Column A (Level) stays "1", same as the cloned row, so only the
# Code and Display columns need to be overwritten.
$wsConcepts.Cells.Item(11, 2).Value = "POLYM"
$wsConcepts.Cells.Item(11, 3).Value = "Polymalformation"
